$d = $word.ActiveDocument

function Replace-InRange($rng, $old, $new) {
    $r = $rng.Duplicate
    $found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
    return $r
}

function Replace-Text($old, $new) {
    Replace-InRange $d.Content $old $new | Out-Null
}

# --- Paragraph 1 (introduction) ---
$old1 = 'Predators and prey have a suite of strategies they use to maximize benefit with the least cost. Predators are able to switch between food sources when one source becomes hard to find (Carle and Rowe 2014). Other strategies can be gleaned from our foraging lab - central place foraging and traplining, in which a predator will center themselves where prey generally congregates or when a predator uses the most efficient path possible. Prey strategies focus decreasing the cost of predation through group foraging, crypticism, or mimicry. With each strategy prey can avoid capture by working as a group, disguising oneself as an inanimate object, or mimicking a toxic species when foraging. Each strategy, either predator or prey, confers some selective pressures on the other groups and on associated species not involved in the interaction.'
$new1 = 'Predators and prey have a suite of strategies they use to maximize benefit with the least cost. Predators are able to switch between food sources when prey is better concealed (Carle and Rowe 2014). Other strategies can be gleaned from our foraging lab - central place foraging and traplining, in which a predator will center themselves where prey generally congregate or when a predator uses the most efficient path possible to capture prey. Prey strategies focus decreasing the cost of predation and optimize foraging success through group foraging, crypticism, or mimicry.  Using these strategies, prey can avoid capture by working as a group, disguising oneself as an inanimate object, or mimicking a toxic species warning coloration. Each strategy, either predator or prey, confers some selective pressures on associated and unassociated species involved in the interaction.'
Replace-Text $old1 $new1

# --- Paragraph 2 (has a lastRenderedPageBreak mid-paragraph; edit each side ---
# --- separately so the empty <w:lastRenderedPageBreak/> run stays untouched) ---
$p2 = $d.Paragraphs.Item(3)
$old2a = 'There is evidence that predators will switch their food source to a toxic prey when undefended prey are better concealed (Carle and Rowe 2014), however the authors stipulate a natural environment could undermine this switch in diet. Their results suggest that this change in diet follows a risk-prone strategy, in which, the predator eats more toxic prey when the probability of finding undefended prey is low. This strategy leads to selective pressure on the predator to tolerate toxins and on the prey to be more cryptic or more toxic. In the foraging lab, simulated predators followed two strategies - central place foraging and traplining - in the first, one predator stayed where prey was most abundant while the latter predator took the most efficient path to capture prey. For central place foraging, this strategy maximizes benefit and minimizes the cost of searching for food. Alternatively, traplining does the same to a lesser degree by maximizing prey capture and reduces the cost to search. '
$new2a = 'There is evidence that predators will switch their food source to a toxic prey when undefended prey are better concealed (Carle and Rowe 2014), however the authors stipulate a natural environment could undermine this switch in diet. Their results suggest that this change in diet follows a risk-prone strategy, in which, the predator eats more toxic prey when the probability of finding undefended prey is low. This strategy leads to selective pressure on the predator to tolerate toxins and on the prey to be more cryptic or more toxic. In the foraging lab, simulated predators followed two strategies - central place foraging and traplining - in the first, one predator stayed where prey was most abundant, while the latter predator took the most efficient path to capture prey. For central place foraging, this strategy maximizes benefit of resource capture and minimizes the cost of searching for food. Alternatively, traplining does something similar by maximizing prey capture and reducing '
$p2BeforeFull = $d.Range($p2.Range.Start, $p2.Range.End)
$p2r = Replace-InRange $p2BeforeFull $old2a $new2a
$p2BreakPos = $p2r.End
$p2now = $d.Paragraphs.Item(3)
$p2After = $d.Range($p2BreakPos, $p2now.Range.End)
$old2b = 'In this method, the predator optimizes their foraging path in a group of prey at the reduced cost to search and increased resource benefit.'
$new2b = 'the cost to search. In this method, the predator optimizes their foraging path in a group of prey at the reduced search cost and increased resource capture benefit.'
Replace-InRange $p2After $old2b $new2b | Out-Null

# --- Paragraph 3 ---
$old3 = 'The foraging lab revealed that group foraging can increase individual benefit at the cost of sharing with group members. Working in a group reduced competition from non-group members because their access to the resource was reduced. Combining efforts increases foraging efficiency for group members by cutting searching time in half and decreases the probability of predation on an individual group member versus a random non-group member. Group size increases visibility, contrary to crypticism where prey disguise themselves as an inanimate object. In rock outcrop communities, invertebrate prey blend in to the background to acquire resources at the least cost (Boggess and Kauffman 2022). Optimal foraging in this way enables prey to adopt a central place foraging strategy to maximize prey capture and predator avoidance, and minimize search time. Mullerian mimicry provides the benefit of reduced predation due to toxicity at the cost of maintaining toxicity and aposemitism. According to Hoyal Cuthill and Charleston (2012), there is evidence that Mullerian mimics co-evolve in a mimicry arms race where the larger population drives divergence and the rarer receives greater benefit. In a predator switching model, mimics forage most optimally when their warning coloration is effective and toxic enough to enforce the warning.'
$new3 = 'From the prey perspective, the foraging lab revealed that group foraging can increase individual benefit at the cost of sharing with group members. Working in a group reduced competition from non-group members because their access to the resource was reduced. Combining efforts increases foraging efficiency for group members by cutting searching time in proportion to the group size and decreasing the probability of predation on an individual group member versus a random non-group member. Group size increases visibility, contrary to crypticism where prey disguise themselves as an inanimate object to reduce visibility. In rock outcrop communities, invertebrate prey blend in to the background to acquire resources at the least cost (Boggess and Kauffman 2022). Optimal foraging in this way enables prey to adopt a central place foraging strategy to maximize prey capture and predator avoidance and minimize search time. Müllerian mimicry provides the benefit of reduced predation due to toxicity at the cost of maintaining toxicity and aposemitism. According to Hoyal Cuthill and Charleston (2012), there is evidence that Müllerian mimics co-evolve in a mimicry arms race where the larger population drives divergence, and the rarer population receives greater benefit. In a predator switching model, mimics forage most optimally when their warning coloration is effective and toxic enough to enforce the warning.'
Replace-Text $old3 $new3

# Append new sentence at end of paragraph 3 (after the bookmarks, before the paragraph mark)
$p3 = $d.Paragraphs.Item(4)
$insPoint = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$append3 = ' This strategy increases resource capture and reduces predation risk on the Müllerian mimics. '
$insPoint.InsertBefore($append3) | Out-Null

# --- Paragraph 4 (conclusion; also has a lastRenderedPageBreak mid-paragraph) ---
$p4 = $d.Paragraphs.Item(5)
$old4a = 'Strategies in which prey can forage optimally through aposemitism, by being cryptic, or in a group at the cost of sharing and benefit of reduced predation. Predator strategies include; traplining, central place foraging, and diet switching - each with respective foraging benefits and costs. The strategies employed by predators and prey provide selection '
$new4a = 'Strategies in which prey can forage optimally include aposemitism, crypticism, or group foraging. While predator strategies include traplining, central place foraging, and diet '
$p4BeforeFull = $d.Range($p4.Range.Start, $p4.Range.End)
$p4r = Replace-InRange $p4BeforeFull $old4a $new4a
$p4BreakPos = $p4r.End
$p4now = $d.Paragraphs.Item(5)
$p4After = $d.Range($p4BreakPos, $p4now.Range.End)
$old4b = 'pressures on effectiveness of predator and prey strategies directly involved in the interaction. Unassociated species experience some selective pressure as a result of reduced predation during switching or increased crypticism or aposemitism, for example. Further research is needed on the indirect effects of co-divergence between predator and prey associated with predation and species unassociated with the interaction. Predator and prey strategies confer benefits and costs to foraging and the interaction between strategies and species could be an interesting avenue for future co-evolutionary research.'
$new4b = 'switching. All of which incur foraging benefits and costs to the predator and the prey. These strategies employed by predators and prey provide selection pressures on effectiveness of predator and prey strategies directly involved in the interaction. Unassociated species experience some selective pressure as a result of reduced predation during switching or increased crypticism or aposemitism, for example. Further research is needed on the indirect effects of co-divergence between predator and prey associated with predation and species unassociated with the interaction. Predator and prey strategies confer benefits and costs to foraging and the interaction between strategies and species could be an interesting avenue for future co-evolutionary research.'
Replace-InRange $p4After $old4b $new4b | Out-Null

Write-Host "Done."
